# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows into the Chirimoya (Vega Modelo de
# Temuco) data table. Both new rows are inserted by shifting the existing
# rows below them down by one, matching how the source system appends each
# week's new observation near the top of the (descending-by-recency) list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper to fill one full data row (columns A..T) at a given row number ---
function Set-ChirimoyaRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value = 10
    $ws.Cells.Item($Row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($Row, 3).Value = "La Araucanía"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 9
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100107
    $ws.Cells.Item($Row, 8).Value = "Otros"
    $ws.Cells.Item($Row, 9).Value = 100107002
    $ws.Cells.Item($Row, 10).Value = "Chirimoya"
    $ws.Cells.Item($Row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# Insert the first new observation at row 111 (pushes the former row 111
# onward down to row 112, etc.)
$ws.Rows.Item(111).Insert()
Set-ChirimoyaRow -Row 111 -Fecha 44846 -Calidad "Primera" -Volumen 50 `
    -PrecioMin 3000 -PrecioMax 3000 -PrecioProm 3000 `
    -Unidad "$/kilo (en caja de 15 kilos)" -Origen "Provincia del Elquí" `
    -PrecioKg 3000 -KgUnidad 1

# Insert the second new observation at row 116 (final numbering) — pushes
# what is currently row 116 onward down to row 117, etc.
$ws.Rows.Item(116).Insert()
Set-ChirimoyaRow -Row 116 -Fecha 44845 -Calidad "Primera" -Volumen 250 `
    -PrecioMin 3000 -PrecioMax 3000 -PrecioProm 3000 `
    -Unidad "$/kilo (en caja de 15 kilos)" -Origen "Provincia del Elquí" `
    -PrecioKg 3000 -KgUnidad 1
